# Re-gen MDS with Stress Monitoring
# Updates the symmetric distance matrix in Sheet1 (rows 2-8, cols B-H)
# with newly computed values, keeping the diagonal at 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New distance matrix values (row label index -> column letter -> value)
$values = @{
    2 = @{ "C" = 0.2013582225080341; "D" = 0.2333163773956504; "E" = 0.4564820952209243; "F" = 0.346880430367577;  "G" = 0.4385702984231664; "H" = 0.3249153431110011 }
    3 = @{ "B" = 0.2013582225080341; "D" = 0.1788672370654924; "E" = 0.429672575826602;  "F" = 0.3275712988552347; "G" = 0.4706396546313344; "H" = 0.3283993633665939 }
    4 = @{ "B" = 0.2333163773956504; "C" = 0.1788672370654924; "E" = 0.464145664093723;  "F" = 0.3209117648019807; "G" = 0.4678895691545143; "H" = 0.3069243545134991 }
    5 = @{ "B" = 0.4564820952209243; "C" = 0.429672575826602;  "D" = 0.464145664093723;  "F" = 0.2590022844278493; "G" = 0.3403533178805759; "H" = 0.3960987825108018 }
    6 = @{ "B" = 0.346880430367577;  "C" = 0.3275712988552347; "D" = 0.3209117648019807; "E" = 0.2590022844278493; "G" = 0.3253573968246197; "H" = 0.2378227310462223 }
    7 = @{ "B" = 0.4385702984231664; "C" = 0.4706396546313344; "D" = 0.4678895691545143; "E" = 0.3403533178805759; "F" = 0.3253573968246197; "H" = 0.3956867380715692 }
    8 = @{ "B" = 0.3249153431110011; "C" = 0.3283993633665939; "D" = 0.3069243545134991; "E" = 0.3960987825108018; "F" = 0.2378227310462223; "G" = 0.3956867380715692 }
}

foreach ($row in $values.Keys) {
    foreach ($col in $values[$row].Keys) {
        $ws.Range("$col$row").Value = $values[$row][$col]
    }
}
